$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 166.74193
$ws.Range("I33").Value = 155.60869
$ws.Range("K33").Value = 155.60869
$ws.Range("M33").Value = 73.39131
$ws.Range("H63").Value = 45271
$ws.Range("J63").Value = 45271
$ws.Range("L63").Value = 45271
$ws.Range("N63").Value = -46519
$ws.Range("H66").Value = 45271
$ws.Range("J66").Value = 45271
$ws.Range("L66").Value = 135813
$ws.Range("N66").Value = -142053
$ws.Range("H129").Value = 1368.9395
$ws.Range("I129").Value = 604.8333
$ws.Range("J129").Value = 1805.5714
$ws.Range("K129").Value = 1814.4999
$ws.Range("L129").Value = 5416.7142
$ws.Range("M129").Value = 3185.5001
$ws.Range("N129").Value = -15416.7142
$ws.Range("H137").Value = 52643604
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 52643604
$ws.Range("K137").Value = 0
$ws.Range("L137").ClearContents()
$ws.Range("M137").Value = 157930812
$ws.Range("N137").Value = -157935912

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 18000
$ws.Range("J34").Value = 18000
$ws.Range("L34").Value = 18000
$ws.Range("N34").Value = -18542
$ws.Range("H61").Value = 5374
$ws.Range("I61").Value = 4514.2856
$ws.Range("J61").Value = 6878.5
$ws.Range("K61").Value = 4514.2856
$ws.Range("L61").Value = 6878.5
$ws.Range("M61").Value = -4302.2856
$ws.Range("N61").Value = -7302.5
$ws.Range("H88").Value = 2966.25
$ws.Range("I88").Value = 1959.6
$ws.Range("J88").Value = 7999.5
$ws.Range("K88").Value = 1959.6
$ws.Range("L88").Value = 7999.5
$ws.Range("M88").Value = -1553.6
$ws.Range("N88").Value = -8811.5
$ws.Range("H91").Value = 2966.25
$ws.Range("I91").Value = 1959.6
$ws.Range("J91").Value = 7999.5
$ws.Range("K91").Value = 1959.6
$ws.Range("L91").Value = 7999.5
$ws.Range("M91").Value = -555.5999999999999
$ws.Range("N91").Value = -10807.5
$ws.Range("H132").Value = 3450
$ws.Range("I132").Value = 3756.5454
$ws.Range("J132").Value = 2968.2856
$ws.Range("K132").Value = 11269.6362
$ws.Range("L132").Value = 8904.856800000001
$ws.Range("M132").Value = -8739.636200000001
$ws.Range("N132").Value = -13964.8568
$ws.Range("H136").Value = 5374
$ws.Range("I136").Value = 4514.2856
$ws.Range("J136").Value = 6878.5
$ws.Range("K136").Value = 13542.8568
$ws.Range("L136").Value = 20635.5
$ws.Range("M136").Value = -10992.8568
$ws.Range("N136").Value = -25735.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H33").Value = 6500
$ws.Range("I33").Value = 3000
$ws.Range("K33").Value = 3000
$ws.Range("M33").Value = -2664
$ws.Range("H80").Value = 726.4211
$ws.Range("I80").Value = 341.2
$ws.Range("J80").Value = 864
$ws.Range("K80").Value = 341.2
$ws.Range("L80").Value = 864
$ws.Range("M80").Value = 656.8
$ws.Range("N80").Value = -2860
$ws.Range("H83").Value = 726.4211
$ws.Range("I83").Value = 341.2
$ws.Range("J83").Value = 864
$ws.Range("K83").Value = 1706
$ws.Range("L83").Value = 4320
$ws.Range("M83").Value = 3286
$ws.Range("N83").Value = -14304
$ws.Range("H94").Value = 913.5172
$ws.Range("I94").Value = 735.5909
$ws.Range("J94").Value = 1472.7142
$ws.Range("K94").Value = 735.5909
$ws.Range("L94").Value = 1472.7142
$ws.Range("M94").Value = -284.5909
$ws.Range("N94").Value = -2374.7142
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").ClearContents()
$ws.Range("N119").Value = 0

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 176178.33
$ws.Range("I4").Value = 1018301
$ws.Range("J4").Value = 7753.8
$ws.Range("K4").Value = 1018301
$ws.Range("L4").Value = 7753.8
$ws.Range("M4").Value = -1018189
$ws.Range("N4").Value = -7977.8
$ws.Range("H31").Value = 2474.347
$ws.Range("I31").Value = 1095.6765
$ws.Range("J31").Value = 5599.3335
$ws.Range("K31").Value = 1095.6765
$ws.Range("L31").Value = 5599.3335
$ws.Range("M31").Value = -800.6765
$ws.Range("N31").Value = -6189.3335
$ws.Range("H34").Value = 2474.347
$ws.Range("I34").Value = 1095.6765
$ws.Range("J34").Value = 5599.3335
$ws.Range("K34").Value = 1095.6765
$ws.Range("L34").Value = 5599.3335
$ws.Range("M34").Value = -893.6765
$ws.Range("N34").Value = -6003.3335

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H127").Value = 2000
$ws.Range("J127").Value = 2000
$ws.Range("L127").Value = 6000
$ws.Range("N127").Value = -15920
$ws.Range("H131").Value = 7152311.5
$ws.Range("J131").Value = 10007093
$ws.Range("L131").Value = 30021279
$ws.Range("N131").Value = -30031359
$ws.Range("H133").Value = 5825.067
$ws.Range("I133").Value = 3240
$ws.Range("J133").Value = 7117.6
$ws.Range("K133").Value = 9720
$ws.Range("L133").Value = 21352.8
$ws.Range("M133").Value = -4660
$ws.Range("N133").Value = -31472.8
$ws.Range("H138").Value = 3140.9
$ws.Range("I138").Value = 3140.9
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 9422.700000000001
$ws.Range("L138").Value = 0
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -4282.700000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 67666.664
$ws.Range("J80").Value = 100000
$ws.Range("L80").Value = 100000
$ws.Range("N80").Value = -101996
$ws.Range("H83").Value = 67666.664
$ws.Range("J83").Value = 100000
$ws.Range("L83").Value = 500000
$ws.Range("N83").Value = -509984
$ws.Range("H132").Value = 3866.6365
$ws.Range("I132").Value = 3509.5
$ws.Range("J132").Value = 4295.2
$ws.Range("K132").Value = 10528.5
$ws.Range("L132").Value = 12885.6
$ws.Range("M132").Value = -7998.5
$ws.Range("N132").Value = -17945.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 338.66666
$ws.Range("I22").Value = 266.66666
$ws.Range("J22").Value = 374.66666
$ws.Range("K22").Value = 266.66666
$ws.Range("L22").Value = 374.66666
$ws.Range("M22").Value = 28.33334000000002
$ws.Range("N22").Value = -964.66666
$ws.Range("H27").Value = 338.66666
$ws.Range("I27").Value = 266.66666
$ws.Range("J27").Value = 374.66666
$ws.Range("K27").Value = 266.66666
$ws.Range("L27").Value = 374.66666
$ws.Range("M27").Value = -159.66666
$ws.Range("N27").Value = -588.66666
$ws.Range("H43").Value = 3691.6667
$ws.Range("J43").Value = 2150
$ws.Range("L43").Value = 2150
$ws.Range("N43").Value = -2536
$ws.Range("H93").Value = 2097.5625
$ws.Range("I93").Value = 1447.65
$ws.Range("J93").Value = 3180.75
$ws.Range("K93").Value = 1447.65
$ws.Range("L93").Value = 3180.75
$ws.Range("M93").Value = -199.6500000000001
$ws.Range("N93").Value = -5676.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4770.636
$ws.Range("I81").Value = 1216.1666
$ws.Range("J81").Value = 9036
$ws.Range("K81").Value = 2432.3332
$ws.Range("L81").Value = 18072
$ws.Range("M81").Value = -1371.3332
$ws.Range("N81").Value = -20194
$ws.Range("H84").Value = 4770.636
$ws.Range("I84").Value = 1216.1666
$ws.Range("J84").Value = 9036
$ws.Range("K84").Value = 12161.666
$ws.Range("L84").Value = 90360
$ws.Range("M84").Value = -6857.666000000001
$ws.Range("N84").Value = -100968
$ws.Range("H106").Value = 30000
$ws.Range("J106").Value = 30000
$ws.Range("L106").Value = 30000
$ws.Range("N106").Value = -32524
